# Auto-applied update: cryptos list refresh (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.972.74"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.745.61"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.77"
$ws.Range("E5").Value = "  +4.71%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5078"
$ws.Range("E7").Value = "  -9.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2752"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06194"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "1.744.06"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07250"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.676"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "25.983.21"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.89"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006863"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "1.970.10"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.389"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.67"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.519"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.76"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.858"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08203"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.649"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.655"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9980"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6172"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.756"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.934"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.76"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3931"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7679"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.003"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.352"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.77"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05345"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.572"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3451"
$ws.Range("E51").Value = "  -1.28%  "
